# Auto update: 2025-12-03 08:54:10
# Update K (최종점수) and N (MACRO_SCORE) columns for rows 2-6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 55.8
$ws.Range("N2").Value = 66.04328690552585

$ws.Range("K3").Value = 50.8
$ws.Range("N3").Value = 66.04328690552585

$ws.Range("K4").Value = 50.8
$ws.Range("N4").Value = 66.04328690552585

$ws.Range("K5").Value = 49.2
$ws.Range("N5").Value = 66.04328690552585

$ws.Range("K6").Value = 48
$ws.Range("N6").Value = 66.04328690552585
